$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells whose new text value looks numeric need to be forced to Text
# format first, otherwise Excel auto-converts them to a floating point number
# (losing the exact original decimal text and introducing FP rounding noise).

$ws.Range("D2").Value = "26.640.46"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "1.598.20"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.51"
$ws.Range("E10").Value = "  -0.86%  "

$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").Value = "1.822.53"
$ws.Range("E12").Value = "  +0.22%  "

$ws.Range("D13").Value = "1.595.30"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.73"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").Value = "26.625.57"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("E18").Value = "  +0.91%  "

$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.42"
$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("E21").Value = "  +4.63%  "

$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("E23").Value = "  -0.56%  "

$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.46"
$ws.Range("E25").Value = "  -0.73%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  -0.25%  "

$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("E30").Value = "  +1.38%  "

$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("E32").Value = "  +0.24%  "

$ws.Range("E33").Value = "  +0.83%  "

$ws.Range("D34").Value = "1.274.57"
$ws.Range("E34").Value = "  -1.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.618"
$ws.Range("E35").Value = "  -7.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("E36").Value = "  +0.61%  "

$ws.Range("E37").Value = "  +0.53%  "

$ws.Range("E38").Value = "  -0.47%  "

$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.07"
$ws.Range("E39").Value = "  +17.98%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.838"
$ws.Range("E40").Value = "  +0.29%  "

$ws.Range("E41").Value = "  +2.36%  "

$ws.Range("E42").Value = "  +0.25%  "

$ws.Range("E43").Value = "  -1.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.01"
$ws.Range("E44").Value = "  +0.68%  "

$ws.Range("D45").Value = "1.735.32"
$ws.Range("E45").Value = "  +0.25%  "

$ws.Range("E46").Value = "  +0.26%  "

$ws.Range("E47").Value = "  -2.48%  "

$ws.Range("E48").Value = "  +3.26%  "

$ws.Range("E49").Value = "  +0.86%  "

$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("E51").Value = "  -1.52%  "
